$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "birotulate isochela"
$ws.Range("A2").Value = "style"

$ws.Range("A2").Select()
